$d = $word.ActiveDocument

# Every "old" string below occurs exactly once in the document, so a
# single (non-"replace all") Find.Execute per pair is unambiguous and
# order-independent.
$pairs = @(
    @("2025-05-29 Thursday", "2025-05-30 Friday"),
    @("11×43=", "33×76="),
    @("43×25=", "69×25="),
    @("75×91=", "96×21="),
    @("43×16=", "94×64="),
    @("72×12=", "36×67="),
    @("50×22=", "43×80="),
    @("75×79=", "42×37="),
    @("87×42=", "41×28="),
    @("97×55=", "28×34="),
    @("11×50=", "25×58="),
    @("53×75=", "28×34="),
    @("51×72=", "33×29="),
    @("37×21=", "67×65="),
    @("47×73=", "18×78="),
    @("78×36=", "70×79="),
    @("81×73=", "11×69="),
    @("48×40=", "29×42="),
    @("75×60=", "65×79="),
    @("77×79=", "23×76="),
    @("74×30=", "61×89="),
    @("91×65=", "79×74="),
    @("68×47=", "18×38="),
    @("69×95=", "48×32="),
    @("94×85=", "44×67="),
    @("96×65=", "27×14=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
